$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 33: replace existing member "陈添楠" with new member "王一鸣"
# (B33 "内容:编写web端用例规约" is left untouched)
$ws.Range("A33").Value = "王一鸣"

# Row 34: keep existing member "邱培松", update the task description
$ws.Range("B34").Value = "内容:设计数据库逻辑模型的entity、物理模型的table"

# Row 35: new member "黄龙强" with his task
$ws.Range("A35").Value = "黄龙强"
$ws.Range("B35").Value = "内容:初步编写android ui代码"

# Row 36: new member "李福森" with his task
$ws.Range("A36").Value = "李福森"
$ws.Range("B36").Value = "内容:初步编写web端普通用户的ui界面代码"
